# Revert "adding term 2.0.0"
$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# --- Include from FSIII sheet ---
$incl = $wb.Worksheets.Item("Include from FSIII")

# Remove the extra concept row (was row 13, the UUID aec684bd-...),
# shifting rows 14-15 up by one.
$incl.Rows.Item(13).Delete()
